$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.889.61'
$ws.Range('E2').Value = '  -1.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.450.44'
$ws.Range('E3').Value = '  -2.92%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.10'
$ws.Range('E5').Value = '  +0.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.28'
$ws.Range('E6').Value = '  -1.29%  '

$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.454.56'
$ws.Range('E9').Value = '  -2.68%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0978'
$ws.Range('E10').Value = '  +0.06%  '

$ws.Range('E11').Value = '  -1.81%  '

$ws.Range('E12').Value = '  -3.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.323'
$ws.Range('E13').Value = '  -2.34%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.888.48'
$ws.Range('E14').Value = '  -2.68%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.849.50'
$ws.Range('E15').Value = '  -1.01%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.72'
$ws.Range('E16').Value = '  -1.77%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000132'
$ws.Range('E17').Value = '  -1.49%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.455.90'
$ws.Range('E18').Value = '  -2.71%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.30'
$ws.Range('E19').Value = '  -3.31%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.11'
$ws.Range('E20').Value = '  -1.18%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '310.72'
$ws.Range('E21').Value = '  -3.49%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('E22').Value = '  -1.10%  '

$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.91'
$ws.Range('E24').Value = '  +0.40%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.403'
$ws.Range('E25').Value = '  -0.94%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.590.46'
$ws.Range('E26').Value = '  -1.37%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.09%  '

$ws.Range('E28').Value = '  -1.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.23'
$ws.Range('E29').Value = '  -2.22%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.44'
$ws.Range('E30').Value = '  +2.83%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0736'
$ws.Range('E31').Value = '  -2.14%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.69'
$ws.Range('E32').Value = '  -1.48%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.22'
$ws.Range('E33').Value = '  -1.21%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.13'
$ws.Range('E34').Value = '  -4.24%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.03%  '

$ws.Range('E36').Value = '  -0.08%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.79'
$ws.Range('E37').Value = '  -1.95%  '

$ws.Range('E38').Value = '  -4.80%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.79'
$ws.Range('E39').Value = '  -3.32%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.814'
$ws.Range('E40').Value = '  +5.83%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.22'
$ws.Range('E41').Value = '  -0.67%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.45'
$ws.Range('E42').Value = '  -2.20%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.40'
$ws.Range('E43').Value = '  -1.47%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '262.37'
$ws.Range('E44').Value = '  -5.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.586'
$ws.Range('E45').Value = '  -1.95%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.81'
$ws.Range('E46').Value = '  -3.47%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0920'
$ws.Range('E47').Value = '  +0.42%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.23'
$ws.Range('E48').Value = '  -5.97%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0493'
$ws.Range('E49').Value = '  -1.29%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0211'
$ws.Range('E50').Value = '  -1.34%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.01'
$ws.Range('E51').Value = '  -3.84%  '
